$wb = $excel.ActiveWorkbook

# --- Sheet1: rename to "Multi User Credentials" ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Multi User Credentials"

# Add new column D (description / testcase) to existing rows.
# "testcase" is written before "description" so the shared-strings table
# ends up with the same insertion ordering as the target workbook.
$ws1.Range("D2").Value = "testcase"
$ws1.Range("D1").Value = "description"
$ws1.Range("D3").Value = "testcase"
$ws1.Range("D4").Value = "testcase"

# Column widths on sheet1
$ws1.Columns.Item(1).ColumnWidth = 9.166666
$ws1.Columns.Item(4).ColumnWidth = 13.166666

# --- Sheet2: add new sheet right after Sheet1 ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Single User Credentials"

$ws2.Range("A1").Value = "testcaseid"
$ws2.Range("B1").Value = "username"
$ws2.Range("C1").Value = "password"
$ws2.Range("D1").Value = "description"

$ws2.Range("A2").Value = 3
$ws2.Range("B2").Value = "standard_user"
$ws2.Range("C2").Value = "secret_sauce"
$ws2.Range("D2").Value = "testcase"

# Give the new sheet's cells the same cell formatting (left/top aligned)
# that the existing data already uses, by copying the format from
# Sheet1 instead of setting alignment properties (which would mint new,
# unused cell styles in styles.xml).
$ws1.Range("A1").Copy() | Out-Null
$ws2.Range("A1:D2").PasteSpecial(-4122) | Out-Null # xlPasteFormats

# Column widths on sheet2
$ws2.Columns.Item(1).ColumnWidth = 8.333333
$ws2.Columns.Item(2).ColumnWidth = 12.666666
$ws2.Columns.Item(3).ColumnWidth = 11.333333
$ws2.Columns.Item(4).ColumnWidth = 11

# Selections on each sheet
$ws1.Range("B21").Select() | Out-Null
$ws2.Range("F22").Select() | Out-Null

# Sheet2 ("Single User Credentials") ends up the active tab
$ws2.Activate()
